# APACHE POI Writing to Excel
# Cell A3 on the "Employees" sheet held "Madam" - correct it to "Adam" and
# move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Adam"
$ws.Range("A2").Select()
